$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple single-run cell text replacements (top of table) ---
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "37"
$t.Cell(6,1).Range.Text  = "0.00051"
$t.Cell(7,1).Range.Text  = "0.00019"

# --- Remove the row that held "0.00023" (row 9) ---
$t.Rows.Item(9).Delete()

# --- After the deletion, rows shift up by one; update the now-shifted rows ---
$t.Cell(10,1).Range.Text = "0.00038"
$t.Cell(11,1).Range.Text = "0.00047"

# --- Insert a new row after row 11 (the one now holding "0.00047") ---
$t.Rows.Add($t.Rows.Item(12)) | Out-Null
$t.Cell(12,1).Range.Text = "0.00854"

# --- Collapse the multi-run summary rows at the bottom into single values ---
$t.Cell(44,1).Range.Text = "99.99"
$t.Cell(45,1).Range.Text = "0.01"
$t.Cell(46,1).Range.Text = "65"
